$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.941.78"
$ws.Range("E2").Value = "'  +0.64%  "
$ws.Range("D3").Value = "'1.746.35"
$ws.Range("E3").Value = "'  -0.26%  "
$ws.Range("D4").Value = "'0.9986"
$ws.Range("E4").Value = "'  -0.24%  "
$ws.Range("D5").Value = "'233.55"
$ws.Range("E5").Value = "'  -1.24%  "
$ws.Range("D6").Value = "'0.9976"
$ws.Range("E6").Value = "'  -0.26%  "
$ws.Range("D7").Value = "'0.5164"
$ws.Range("E7").Value = "'  +1.94%  "
$ws.Range("D8").Value = "'0.2811"
$ws.Range("E8").Value = "'  +7.45%  "
$ws.Range("D9").Value = "'39.65"
$ws.Range("D10").Value = "'0.06125"
$ws.Range("E10").Value = "'  -1.12%  "
$ws.Range("D11").Value = "'1.739.90"
$ws.Range("D12").Value = "'0.07019"
$ws.Range("E12").Value = "'  +1.27%  "
$ws.Range("D13").Value = "'15.44"
$ws.Range("E13").Value = "'  +0.04%  "
$ws.Range("D14").Value = "'0.6406"
$ws.Range("E14").Value = "'  +5.77%  "
$ws.Range("D15").Value = "'4.521"
$ws.Range("E15").Value = "'  +1.33%  "
$ws.Range("D16").Value = "'76.94"
$ws.Range("E16").Value = "'  -2.39%  "
$ws.Range("D17").Value = "'0.9969"
$ws.Range("E17").Value = "'  -0.35%  "
$ws.Range("D18").Value = "'0.9967"
$ws.Range("E18").Value = "'  -0.34%  "
$ws.Range("D19").Value = "'25.956.17"
$ws.Range("E19").Value = "'  +0.56%  "
$ws.Range("D20").Value = "'11.50"
$ws.Range("E20").Value = "'  -1.59%  "
$ws.Range("D21").Value = "'0.000006619"
$ws.Range("E21").Value = "'  -1.63%  "
$ws.Range("D22").Value = "'1.963.01"
$ws.Range("D23").Value = "'4.140"
$ws.Range("E23").Value = "'  +2.14%  "
$ws.Range("D24").Value = "'8.551"
$ws.Range("E24").Value = "'  +4.26%  "
$ws.Range("D25").Value = "'5.152"
$ws.Range("E25").Value = "'  -0.06%  "
$ws.Range("D27").Value = "'1.506"
$ws.Range("E27").Value = "'  +3.57%  "
$ws.Range("D28").Value = "'1.834"
$ws.Range("E28").Value = "'  +1.48%  "
$ws.Range("D29").Value = "'15.08"
$ws.Range("E29").Value = "'  -0.18%  "
$ws.Range("D30").Value = "'103.16"
$ws.Range("E30").Value = "'  +0.92%  "
$ws.Range("D31").Value = "'0.08314"
$ws.Range("E31").Value = "'  +0.57%  "
$ws.Range("D32").Value = "'3.639"
$ws.Range("E32").Value = "'  -1.79%  "
$ws.Range("D33").Value = "'3.433"
$ws.Range("E33").Value = "'  +1.03%  "
$ws.Range("D34").Value = "'0.04414"
$ws.Range("E34").Value = "'  +1.44%  "
$ws.Range("D35").Value = "'2.614"
$ws.Range("E35").Value = "'  -1.23%  "
$ws.Range("D36").Value = "'0.9838"
$ws.Range("E36").Value = "'  -1.63%  "
$ws.Range("D37").Value = "'0.6099"
$ws.Range("E37").Value = "'  +1.44%  "
$ws.Range("D38").Value = "'2.684"
$ws.Range("E38").Value = "'  -0.68%  "
$ws.Range("D39").Value = "'0.01572"
$ws.Range("E39").Value = "'  +1.46%  "
$ws.Range("D40").Value = "'1.935"
$ws.Range("E40").Value = "'  -1.35%  "
$ws.Range("D41").Value = "'0.9963"
$ws.Range("E41").Value = "'  -0.37%  "
$ws.Range("D42").Value = "'100.65"
$ws.Range("E42").Value = "'  -2.36%  "
$ws.Range("D43").Value = "'0.3869"
$ws.Range("E43").Value = "'  +1.47%  "
$ws.Range("D44").Value = "'0.7351"
$ws.Range("E44").Value = "'  -2.96%  "
$ws.Range("D45").Value = "'4.973"
$ws.Range("E45").Value = "'  +2.49%  "
$ws.Range("D46").Value = "'0.05454"
$ws.Range("E46").Value = "'  -0.54%  "
$ws.Range("D47").Value = "'6.368"
$ws.Range("E47").Value = "'  +7.29%  "
$ws.Range("D48").Value = "'0.1118"
$ws.Range("E48").Value = "'  +3.73%  "
$ws.Range("D49").Value = "'52.71"
$ws.Range("E49").Value = "'  +1.43%  "
$ws.Range("D50").Value = "'29.87"
$ws.Range("E50").Value = "'  -0.80%  "
$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.549"
$ws.Range("E51").Value = "'  +1.50%  "
